# fallas en agregado a excel
# Add two more scored responses (rows 3 and 4) to the Zulliger scoring sheet,
# and fix the "Cont" value that was missing in row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix row 2 (Lam 1): Cont column was missing "Ad," ---
$ws.Range("I2").Value = " A Ad,"

# --- Row 3 (Lam 2) ---
$ws.Range("A3").Value = 2
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "1"
$ws.Range("C3").Value = "?"
$ws.Range("D3").Value = "?"
$ws.Range("E3").Value = "?"
$ws.Range("F3").Value = " m, C"
$ws.Range("G3").Value = "?"
$ws.Range("H3").Value = "?"
$ws.Range("I3").Value = " Fi Hx,"
$ws.Range("J3").Value = "?"
$ws.Range("K3").Value = "?"
$ws.Range("L3").Value = "?"

# --- Row 4 (Lam 3) ---
$ws.Range("A4").Value = 3
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "1"
$ws.Range("C4").Value = "?"
$ws.Range("D4").Value = "?"
$ws.Range("E4").Value = "?"
$ws.Range("F4").Value = " M, M"
$ws.Range("G4").Value = "?"
$ws.Range("H4").Value = "?"
$ws.Range("I4").Value = " H H, Hd,"
$ws.Range("J4").Value = "Po3"
$ws.Range("K4").Value = "?"
$ws.Range("L4").Value = "?"
